# Add new database-structure rows to the "Tables" worksheet describing
# additional columns for the "Student" table plus new "Class" and "School"
# tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")

# Data to append, starting at row 14.
# Columns: A=Table ID, B=TableName, C=ColumnName, D=DataType, E=Size, F=Unique Key, G=Foreign Key
$rows = @(
    @(2, "Student", "UIDAI number", "Int",      $null, $null, $null),
    @(2, "Student", "Address",      "nvarchar", 600,   $null, $null),
    @(3, "Class",   "ID",           "Int",      $null, $null, $null),
    @(3, "Class",   "Class_Number", "nvarchar", 200,   $null, $null),
    @(4, "School",  "ID",           "Int",      $null, $null, $null),
    @(4, "School",  "Name",         "nvarchar", 200,   $null, $null),
    @(4, "School",  "Address",      "nvarchar", 600,   $null, $null),
    @(4, "School",  "Principal",    "Int",      $null, $null, "User.UserID"),
    @(4, "School",  "RU",           "Int",      $null, $null, $null),
    @(4, "School",  "District",     "nvarchar", 200,   $null, $null)
)

$startRow = 14
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    if ($null -ne $data[4]) {
        $ws.Cells.Item($r, 5).Value = $data[4]
    }
    if ($null -ne $data[5]) {
        $ws.Cells.Item($r, 6).Value = $data[5]
    }
    if ($null -ne $data[6]) {
        $ws.Cells.Item($r, 7).Value = $data[6]
    }
}

$ws.Range("A24").Select()
